$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.258.38"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.882.47"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'243.68"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.4903"
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").Value = "'0.2917"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "'0.06616"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "1.880.90"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'16.46"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").Value = "'0.07222"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'0.6679"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'4.990"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "30.169.24"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "'0.000007821"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'12.81"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "2.116.72"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'0.9974"
$ws.Range("D22").Value = "'4.768"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "'5.889"
$ws.Range("E23").Value = "  +5.95%  "
$ws.Range("D24").Value = "'9.206"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'151.72"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").Value = "'143.55"
$ws.Range("E26").Value = "  +6.16%  "
$ws.Range("D27").Value = "'17.01"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").Value = "'1.897"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "'1.401"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "'4.222"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'0.08821"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "'3.978"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").Value = "'0.05212"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "'0.7251"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").Value = "'1.116"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "'2.662"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  +12.23%  "
$ws.Range("D38").Value = "'2.683"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'2.178"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "'0.9341"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'0.4270"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").Value = "'104.42"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'5.777"
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'7.463"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "'0.1287"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "'0.05736"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'32.89"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "'8.289"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "'0.3787"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").Value = "'1.355"
$ws.Range("E51").Value = "  +1.60%  "
